# Updated cryptos list on Fri Jun 30 04:39:40 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row, and
# swap the Avalanche/Dai rows (18/19) back into rank order.
# Price cells whose new text would otherwise be auto-parsed as a number by
# Excel (e.g. "0.9974", "1.000") are forced to Text format first so the
# stored value keeps its original string form (leading/trailing zeros etc.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.840.71'
$ws.Range("E2").Value = '  +2.17%  '
$ws.Range("D3").Value = '1.883.75'
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9974'
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.65'
$ws.Range("E5").Value = '  +2.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9987'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4794'
$ws.Range("E7").Value = '  +2.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2841'
$ws.Range("E8").Value = '  +5.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06523'
$ws.Range("E9").Value = '  +4.10%  '
$ws.Range("E10").Value = '  +17.24%  '
$ws.Range("D11").Value = '1.870.08'
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07519'
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.33'
$ws.Range("E13").Value = '  +14.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.109'
$ws.Range("E14").Value = '  +3.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6539'
$ws.Range("E15").Value = '  +5.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '297.31'
$ws.Range("E16").Value = '  +31.28%  '
$ws.Range("D17").Value = '30.765.84'
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.14'
$ws.Range("E18").Value = '  +6.75%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9994'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007562'
$ws.Range("E20").Value = '  +3.91%  '
$ws.Range("D21").Value = '2.116.52'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.168'
$ws.Range("E23").Value = '  +6.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.120'
$ws.Range("E24").Value = '  +4.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '169.08'
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.256'
$ws.Range("E26").Value = '  +0.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.58'
$ws.Range("E27").Value = '  +10.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.968'
$ws.Range("E28").Value = '  +4.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1055'
$ws.Range("E29").Value = '  +1.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.366'
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.143'
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.951'
$ws.Range("E32").Value = '  +4.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04992'
$ws.Range("E33").Value = '  +3.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.178'
$ws.Range("E34").Value = '  +3.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7229'
$ws.Range("E35").Value = '  +2.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.701'
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01932'
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.719'
$ws.Range("E38").Value = '  +2.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.070'
$ws.Range("E39").Value = '  +7.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8911'
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '107.60'
$ws.Range("E41").Value = '  +3.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9997'
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4212'
$ws.Range("E43").Value = '  +5.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.598'
$ws.Range("E44").Value = '  +1.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.354'
$ws.Range("E45").Value = '  +5.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.47'
$ws.Range("E46").Value = '  +9.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1231'
$ws.Range("E47").Value = '  +3.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.65'
$ws.Range("E48").Value = '  +5.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.853'
$ws.Range("E49").Value = '  +3.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.391'
$ws.Range("E50").Value = '  +2.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05564'
$ws.Range("E51").Value = '  +1.05%  '
